$wb = $excel.ActiveWorkbook

# --- Settings sheet ---
$wsSettings = $wb.Worksheets.Item("Settings")
$wsSettings.Range("B1").Value = 5
$wsSettings.Range("B4").Value = 0

# --- Settings_recources sheet ---
$wsRecources = $wb.Worksheets.Item("Settings_recources")
$wsRecources.Range("B28").Value = "xlsm"
# B30 holds a numeric-looking value but must stay text (like "80" did),
# so force text format before assigning to avoid Excel's automatic
# number coercion.
$wsRecources.Range("B30").NumberFormat = "@"
$wsRecources.Range("B30").Value = "60"
$wsRecources.Range("B31").Value = "ne"
